$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.139.98"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.74%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.632.14"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.77%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.00"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.41%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -4.23%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.631.48"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.81%  "
$ws.Range("E10").Value = "  -1.66%  "
$ws.Range("E11").Value = "  +1.36%  "
$ws.Range("E12").Value = "  -0.60%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.24"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.60"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.71%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.110.72"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.86%  "
$ws.Range("E16").Value = "  -2.72%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.009.50"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.86%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.596.84"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.12"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.80%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.04"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.73%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "359.91"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.35"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.84%  "
$ws.Range("E23").Value = "  -5.25%  "
$ws.Range("E24").Value = "  +9.44%  "
$ws.Range("E25").Value = "  -5.73%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("E27").Value = "  -3.08%  "
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.32%  "
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0000101"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.50%  "
$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "554.34"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.13%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.95"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.72%  "
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.37"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.40%  "
$ws.Range("B34").Value = "PancakeSwap"
$ws.Range("C34").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.90"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.54%  "
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.134"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.99%  "
$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.51"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.81%  "
$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "157.49"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.37%  "
$ws.Range("B39").Value = "EthereumClassic"
$ws.Range("C39").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.20"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.13%  "
$ws.Range("B40").Value = "PolygonEcosystemToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.367"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.47%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.22"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.71%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.80"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.64%  "
$ws.Range("B43").Value = "WhiteBITCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.91"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.43%  "
$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.47"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.00%  "
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.16"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.41%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₆0300"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.28%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.588"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.23%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "152.17"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.39%  "
$ws.Range("B50").Value = "Filecoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.81"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.94%  "
$ws.Range("B51").Value = "Optimism"
$ws.Range("C51").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.72"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.34%  "

Write-Output "done"